# Insert a new data row into the "Hortaliza, Feria Lagunitas de Puerto Montt -
# Pepino ensalada" sheet. The new record is inserted at row 114 (pushing the
# existing rows 114-216 down to 115-217, and extending the used range from
# A1:R216 to A1:R217).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 114:216 down by one row to make room for the new record.
$ws.Rows("114:114").Insert()

# Populate the newly inserted row with the new price report.
$ws.Cells.Item(114, 1).Value  = 4
$ws.Cells.Item(114, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(114, 3).Value  = "Los Lagos"
$ws.Cells.Item(114, 4).Value  = 44589
$ws.Cells.Item(114, 5).Value  = 10
$ws.Cells.Item(114, 6).Value  = 100112043
$ws.Cells.Item(114, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(114, 8).Value  = "Sin especificar"
$ws.Cells.Item(114, 9).Value  = "Primera"
$ws.Cells.Item(114, 10).Value = 400
$ws.Cells.Item(114, 11).Value = 14000
$ws.Cells.Item(114, 12).Value = 14000
$ws.Cells.Item(114, 13).Value = 14000
$ws.Cells.Item(114, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(114, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(114, 16).Value = 233
$ws.Cells.Item(114, 17).Value = 60
$ws.Cells.Item(114, 18).Value = "Hortaliza"
